# owner pdf link send to the client
# Update the "OFF ROAD DAY" (column L) figures for the affected vehicle
# rows, then leave the sheet scrolled/selected where the user last left
# off (bottom of the sheet, cell U39) exactly like the authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$offRoadDays = @{
    2  = 1
    4  = 2
    5  = 4
    9  = 2
    10 = 4
    11 = 5
    12 = 1
    13 = 2
    17 = 2
    18 = 4
    21 = 5
    24 = 6
    26 = 3
    28 = 1
    31 = 1
    35 = 1
    38 = 2
    40 = 3
}

foreach ($row in $offRoadDays.Keys) {
    $ws.Range("L$row").Value = $offRoadDays[$row]
}

# Match the author's final view state: scrolled down to row 34 and
# the active selection on U39.
$null = $ws.Range("U39").Select()
